# Sprint Stories status update: implement feature to view previous month expenses.
# - Row 9  ("I wan to send bills to other housemates"): IN PROGRESS -> DONE
# - Row 10 ("I want to copy over existing expenses to new month"): NOT STARTED -> IN PROGRESS

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 should pick up the "IN PROGRESS" look that currently lives on C9,
# so grab that formatting first (copy/paste formats), then overwrite values.
$ws.Range("C9").Copy()
$ws.Range("C10").PasteSpecial(-4122)

# Row 9 should pick up the "DONE" look that already exists on C8.
$ws.Range("C8").Copy()
$ws.Range("C9").PasteSpecial(-4122)

$ws.Range("C9").Value = "DONE"
$ws.Range("C10").Value = "IN PROGRESS"

$excel.CutCopyMode = $false

$ws.Range("F7").Select()
